$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 09:55"

# Row 6 - India
$ws.Range("B6").Value = 908258
$ws.Range("C6").Value = 613
$ws.Range("D6").Value = 572280
$ws.Range("E6").Value = 312242
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 23736

# Row 7 - Rusia
$ws.Range("B7").Value = 739947
$ws.Range("C7").Value = 6248
$ws.Range("D7").Value = 512825
$ws.Range("E7").Value = 215508
$ws.Range("G7").Value = 175
$ws.Range("H7").Value = 11614

# Row 33 - Belgica
$ws.Range("B33").Value = 62781
$ws.Range("C33").Value = 74
$ws.Range("D33").Value = 17223
$ws.Range("E33").Value = 35771
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 9787

# Row 53 - Armenia
$ws.Range("B53").Value = 32490
$ws.Range("C53").Value = 339
$ws.Range("D53").Value = 20729
$ws.Range("E53").Value = 11180
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 581

# Row 90 - Estado de Palestina
$ws.Range("B90").Value = 6579
$ws.Range("C90").Value = 13
$ws.Range("E90").Value = 5453
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 42

# Row 99 - Hungria
$ws.Range("B99").Value = 4258
$ws.Range("C99").Value = 11
$ws.Range("D99").Value = 3106
$ws.Range("E99").Value = 557

# Row 117 - Estonia
$ws.Range("B117").Value = 2015
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 1897
$ws.Range("E117").Value = 49

# Row 124 - Cabo Verde
$ws.Range("B124").Value = 1722
$ws.Range("C124").Value = 24
$ws.Range("D124").Value = 772

# Row 137 - Letonia
$ws.Range("D137").Value = 1022
$ws.Range("E137").Value = 121
